$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("education")

# Rename the "description" header (C1) to "descriptions"
$ws.Range("C1").Value = "descriptions"

# Replace the Lorem-ipsum placeholder in C2 with the JSON-style bullet list
$descriptions = "[`n""Created a board game to simulate skills without coding"",`n""Created story boards to illustrate the user's experience"",`n""Designed and implemented virtual worlds using the Unity game engine"",`n""Developed believable physics, partical effects, reward system, and state machine logic"",`n""Learned about history, psychology, math, mythology, and science in relation to Game Design""`n]"
$ws.Range("C2").Value = $descriptions

# Turn on word-wrap for the whole descriptions column (header + data)
$ws.Range("C1:C2").WrapText = $true

# Match the author's column width / row heights from the diff
$ws.Range("C1").ColumnWidth = 13.7109375
$ws.Range("A1:K1").RowHeight = 15
$ws.Range("A2:K2").RowHeight = 15

# Restore the selection / view state touched in the diff
$ws.Range("C3").Select()
